{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text,items/style\");\nawait context.sync();\n\nconst oldText = \"Jag arbetade med att utveckla servern mottog, lagrade och analyserade registrerad data under anv\u00e4ndarnas cykelturer. Teknikerna var node/javascript, MongoDb, git och testramverket Jest.\";\nconst firstText = \"Jag arbetade med att utveckla en server som tog emot, lagrade och analyserade registrerad data under anv\u00e4ndarnas cykelturer.\";\nconst secondText = \"Teknikerna var node/javascript, MongoDb, git och testramverket Jest.\";\n\nlet target = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === oldText) {\n    target = paras.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Target paragraph not found\");\n}\n\n// Replace the paragraph's text with the first sentence.\ntarget.insertText(firstText, \"Replace\");\n\n// Insert a new paragraph right after it; it inherits the source\n// paragraph's \"ListBullet\" style, carrying the second sentence.\ntarget.insertParagraph(secondText, \"After\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$oldText = \"Jag arbetade med att utveckla servern mottog, lagrade och analyserade registrerad data under anv\u00e4ndarnas cykelturer. Teknikerna var node/javascript, MongoDb, git och testramverket Jest.\"\n$newFirst = \"Jag arbetade med att utveckla en server som tog emot, lagrade och analyserade registrerad data under anv\u00e4ndarnas cykelturer.\"\n$newSecond = \"Teknikerna var node/javascript, MongoDb, git och testramverket Jest.\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.MatchWildcards = $false\n\nif ($find.Execute()) {\n    # Replacing the whole matched range (which spans the paragraph's text,\n    # not its end-of-paragraph mark) with text containing an embedded\n    # carriage return splits it into two paragraphs. Both inherit the\n    # original paragraph's \"ListBullet\" style.\n    $rng.Text = $newFirst + \"`r\" + $newSecond\n}\n"}
